$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.999391737847945, -168212921.728772, 0.922159671674992, 0.999391315109391, -0.922167519616221, 110896.139792307, 7064324.67494998, -2640636.63616957, 0.924928169982639, 0.999548213774337, -0.922167519616221, 46.1093809472369, 2928.93550665638, -2640636.63616957),
    @(0.960692208521904, 191.766895909241, -0.975261734804538, -0.962466664394755, 0.96849110342503, -11.7377292166049, -680.884385471807, 3.2889723996093, -0.972686823897195, -0.83308458247653, 0.96849110342503, -0.0481424735795197, -2.42364630721163, 3.2889723996093),
    @(-0.999391742941703, -168213116.810482, 0.922167574618194, 0.999391739254949, -0.922167621188412, 110907.217803063, 7064972.81960954, -2640639.9759917, 0.922158487999296, 0.999391020084426, -0.922167621188412, 45.1845437846978, 2878.35630414182, -2640639.9759917),
    @(0.961809438901286, 194.081710310473, -0.97466736170898, -0.964998854783382, 0.968248888269855, -12.0780107556611, -702.895176981353, 3.32397544008159, -0.99689925165445, -0.941231183566587, 0.968248888269855, -0.0751628374609167, -4.17131493677733, 3.32397544008159)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $row[$j]
    }
}
